$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collect the full set of cells being updated, covering D2:E51
$updateRange = $ws.Range("D2:E51")

# Force Text format so values like "28.412.04" / "0.4715" are stored as strings,
# matching the original inline-string cell content (not converted to numbers).
$updateRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.412.04"
$ws.Range("E2").Value = "  +3.36%  "
$ws.Range("D3").Value = "1.870.49"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "339.13"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4715"
$ws.Range("E7").Value = "  +2.07%  "
$ws.Range("D8").Value = "0.3979"
$ws.Range("E8").Value = "  +3.83%  "
$ws.Range("D9").Value = "47.77"
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("D10").Value = "0.08038"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("E11").Value = "  +2.80%  "
$ws.Range("E12").Value = "  +4.32%  "
$ws.Range("D13").Value = "1.884.15"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "6.040"
$ws.Range("E14").Value = "  +2.37%  "
$ws.Range("D15").Value = "7.272"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").Value = "91.12"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "0.00001041"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "0.06635"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "17.62"
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "28.419.48"
$ws.Range("E22").Value = "  +3.43%  "
$ws.Range("D23").Value = "5.481"
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("D24").Value = "11.07"
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("E25").Value = "  -2.09%  "
$ws.Range("D26").Value = "2.105.39"
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("D27").Value = "160.69"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("D28").Value = "19.80"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").Value = "2.124"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").Value = "5.505"
$ws.Range("E30").Value = "  +3.70%  "
$ws.Range("D31").Value = "120.33"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "0.9768"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("D33").Value = "0.09525"
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  +4.72%  "
$ws.Range("D36").Value = "5.356"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").Value = "0.06110"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("D38").Value = "0.02252"
$ws.Range("E38").Value = "  +2.55%  "
$ws.Range("D39").Value = "8.370"
$ws.Range("E39").Value = "  +3.72%  "
$ws.Range("D40").Value = "1.184"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "0.1883"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "10.35"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("D46").Value = "0.5593"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").Value = "12.14"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  +4.32%  "
$ws.Range("D49").Value = "0.06954"
$ws.Range("E49").Value = "  +4.56%  "
$ws.Range("D50").Value = "2.071"
$ws.Range("E50").Value = "  +15.80%  "
$ws.Range("D51").Value = "111.93"
$ws.Range("E51").Value = "  +1.28%  "

# Restore the default (Normal) cell style so no stray number-format styling remains
$updateRange.Style = "Normal"
